$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 439, pushing the existing rows 439-452
# down to 441-454 (values unchanged, only their row numbers shift).
$ws.Rows.Item(439).Insert()
$ws.Rows.Item(439).Insert()

# Row 439: new "Primera" quality record, dated 45041
$ws.Cells.Item(439, 1).Value = 1
$ws.Cells.Item(439, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(439, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(439, 4).Value = 45041
$ws.Cells.Item(439, 5).Value = 15
$ws.Cells.Item(439, 6).Value = 100114014
$ws.Cells.Item(439, 7).Value = "Betarraga"
$ws.Cells.Item(439, 8).Value = "Sin especificar"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 800
$ws.Cells.Item(439, 11).Value = 600
$ws.Cells.Item(439, 12).Value = 700
$ws.Cells.Item(439, 13).Value = 650
$ws.Cells.Item(439, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(439, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(439, 16).Value = 162
$ws.Cells.Item(439, 17).Value = 4
$ws.Cells.Item(439, 18).Value = "Hortaliza"

# Row 440: new "Segunda" quality record, dated 45041
$ws.Cells.Item(440, 1).Value = 1
$ws.Cells.Item(440, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(440, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(440, 4).Value = 45041
$ws.Cells.Item(440, 5).Value = 15
$ws.Cells.Item(440, 6).Value = 100114014
$ws.Cells.Item(440, 7).Value = "Betarraga"
$ws.Cells.Item(440, 8).Value = "Sin especificar"
$ws.Cells.Item(440, 9).Value = "Segunda"
$ws.Cells.Item(440, 10).Value = 900
$ws.Cells.Item(440, 11).Value = 600
$ws.Cells.Item(440, 12).Value = 700
$ws.Cells.Item(440, 13).Value = 650
$ws.Cells.Item(440, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(440, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(440, 16).Value = 130
$ws.Cells.Item(440, 17).Value = 5
$ws.Cells.Item(440, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date style (s="2") as the rest
# of column D; Rows.Insert already propagates it from the row above, but
# set it explicitly to be safe.
$ws.Range("D439:D440").NumberFormat = "YYYY-MM-DD HH:MM:SS"
